$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.252.80'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.39%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.252.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.43%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.09%  '

$ws.Range("E7").Value = '  -0.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.243.77'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.54%  '

$ws.Range("E9").Value = '  +6.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.00%  '

$ws.Range("E11").Value = '  +7.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.488'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.74%  '

$ws.Range("E14").Value = '  +7.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.774.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.317.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '553.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +14.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.251.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.17%  '

$ws.Range("E19").Value = '  +3.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.745'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +11.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +20.24%  '

$ws.Range("E28").Value = '  +9.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.26%  '

$ws.Range("E31").Value = '  +7.18%  '

$ws.Range("E32").Value = '  -0.28%  '

$ws.Range("E33").Value = '  +6.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '563.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.28%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '55.25'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0454'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0865'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.131'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.01'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.211.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +11.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.54%  '

$ws.Range("E44").Value = '  +16.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.36%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₃0560'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.35%  '

$ws.Range("E50").Value = '  +4.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.55%  '
